$d = $word.ActiveDocument

function Set-ParagraphText {
    param($paragraph, [string]$newText)
    $r = $paragraph.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Delete()
    $r.InsertBefore($newText)
}

# Paragraphs 1-7 (1-indexed): in-place text swap (keeps paragraph mark/props intact)
Set-ParagraphText $d.Paragraphs.Item(1) "⚡️🚀המאמר היומי של מייק 20.06.24:⚡️🚀"
Set-ParagraphText $d.Paragraphs.Item(2) "WARM: On the Benefits of Weight Averaged Reward Models"
Set-ParagraphText $d.Paragraphs.Item(3) "הסקירה הזו ממשיכה את קו הסקירות בנושא שיפור ביצועי RLHF לטיוב מודלי שפה. כבר דיברנו בסקירות הקודמות על כך שבמהלך RLHF המודל יכול לבצע reward hacking כלומר להתכנס לפוליסי (משקלי המודל) שממקסם את ה-reward ובאותו הזמן יוצר תשובות באיכות ירודה לפרומפטים. "
Set-ParagraphText $d.Paragraphs.Item(4) "המאמר שנסקור קצרות היום מציע לאמן כמה מודלי reward שונים ולהשתמש בממוצע שלהם כ-reward יותר ״יציב״ שעשוי למנוע מהמודל לעשות reward hacking. הבעיה העיקרית בגישה הזאת נובעת מכך שהיא מצריכה להחזיק בזמן אימון RLHF כמה מודלי reward שכמובן דורש יותר משאבי חישוב (ומייקר את חשבון החשמל). "
Set-ParagraphText $d.Paragraphs.Item(5) "המחברים מציע לשלב את התוצאה של המודלים אלא הביצועים שלהם. בשפה פשוטה הם מאמנים כמה מודלי reward וממצעים את המשקלים שלהם. זה מסתמך על איזושהי תופעה שלא ידעתי עליה שנקראת ״Linear mode connectivity״ או LMC הטוענת שהביצועים של מודל עם סכום ממושקל של המשקלים של כמה מודלים אחרים הוא יותר טוב מסכום ממושקל (עם אותם משקלים) של ביצועי המודלים (אולי אתעמק בזה בהמשך). "
Set-ParagraphText $d.Paragraphs.Item(6) "עכשיו כדי לבצע את הפעולה הזו הרשתות צריכות להיות בעלי אותה ארכיטקטורה ומה שונה בין מודלי reward כאן הם פרמטרי אימון כמו קצב למידה ודרופאאוט, סדר שונה של הכנסת דאטה לאימון (סיד שונה כנראה) וגם איתחולים שונים (לוקחים מודלים אחרי צ'קפוינטים שונים ב-SFT)."
Set-ParagraphText $d.Paragraphs.Item(7) "כתוצאה מקבלים מודל reward אחד טוב יותר שמשמש אותם לאימון RLHF."

# Remove the five obsolete paragraphs (old items 8-12) entirely, including their marks
$pStart = $d.Paragraphs.Item(8)
$pEnd = $d.Paragraphs.Item(12)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete() | Out-Null

# Update the final (link) paragraph, which is now item 8
Set-ParagraphText $d.Paragraphs.Item(8) "https://arxiv.org/abs/2401.12187"

Write-Output $d.Paragraphs.Count
